$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 26.51420687852702; "C2" = 19.97425191547378; "D2" = 6.008491069640167; "E2" = 16.31158006290931; "G2" = 3.673851184767146; "N2" = 21.15869419480933
    "B3" = 25.53506499953264; "C3" = 18.99652904208748; "D3" = 5.898023069255254; "E3" = 15.38411371669625; "G3" = 3.681640217266894; "N3" = 21.02631633553118
    "B4" = 24.92636596317875; "C4" = 18.37898605555558; "D4" = 5.831572695935431; "E4" = 14.7932007304306; "G4" = 3.686638508909442; "N4" = 20.94588154569605
    "B5" = 24.67685201298982; "C5" = 18.12337967830436; "D5" = 5.804874021953246; "E5" = 14.54727872259114; "G5" = 3.688730027775257; "N5" = 20.91332211841607
    "B6" = 24.63534500316673; "C6" = 18.08070981099538; "D6" = 5.800464683044735; "E6" = 14.5061433297505; "G6" = 3.689080637043328; "N6" = 20.90792899872004
    "B7" = 24.9230062753449; "C7" = 18.37555432376747; "D7" = 5.831211043613751; "E7" = 14.78990447253269; "G7" = 3.686666493946291; "N7" = 20.94544154637783
    "B8" = 26.17839815314228; "C8" = 19.6409157296671; "D8" = 5.970136406583868; "E8" = 15.99638800390183; "G8" = 3.676492304710149; "N8" = 21.11287833774615
    "B9" = 28.56398378101553; "C9" = 21.97180521676237; "D9" = 6.251984177019161; "E9" = 18.21770676935669; "G9" = 3.658234015458674; "N9" = 21.44766882651193
    "B10" = 30.25045632413138; "C10" = 23.57736003608181; "D10" = 6.462778638751399; "E10" = 19.88784415419331; "G10" = 3.645825020935097; "N10" = 21.69720522513668
    "B11" = 30.99985928464922; "C11" = 24.28216669849652; "D11" = 6.559066181512739; "E11" = 20.60800101663241; "G11" = 3.640392225937446; "N11" = 21.81141568462883
    "B12" = 31.28084151399185; "C12" = 24.54522954660797; "D12" = 6.595551420759634; "E12" = 20.87506636176481; "G12" = 3.63836498885323; "N12" = 21.8547561901627
    "B13" = 31.2204550461224; "C13" = 24.48874685325677; "D13" = 6.587693159327261; "E13" = 20.81779918334653; "G13" = 3.638800261787373; "N13" = 21.84541810520479
    "B14" = 31.02303326570442; "C14" = 24.30388662427506; "D14" = 6.562067566261404; "E14" = 20.63008542853018; "G14" = 3.640224844445886; "N14" = 21.81497954997606
    "B15" = 30.9017353567615; "C15" = 24.1901512324367; "D15" = 6.54637320102055; "E15" = 20.51437206199922; "G15" = 3.641101341381816; "N15" = 21.79634671465571
    "B16" = 30.20110171166966; "C16" = 23.53077078611603; "D16" = 6.456491076374125; "E16" = 19.83998781434145; "G16" = 3.64618429629908; "N16" = 21.68975457115377
    "B17" = 29.76654771189341; "C17" = 23.1195932640692; "D17" = 6.401428312661881; "E17" = 19.41616972614779; "G17" = 3.649356530245039; "N17" = 21.62453582145292
    "B18" = 29.51494243612995; "C18" = 22.88069544491675; "D18" = 6.369796884021234; "E18" = 19.16867060173518; "G18" = 3.65120110887022; "N18" = 21.58708967111163
    "B19" = 29.42947583355726; "C19" = 22.79940188638683; "D19" = 6.35909481200681; "E19" = 19.08422942702175; "G19" = 3.651829099216403; "N19" = 21.57442262930251
    "B20" = 29.81298053241309; "C20" = 23.16361336157917; "D20" = 6.407286027942141; "E20" = 19.46167157136414; "G20" = 3.649016774556777; "N20" = 21.63147170383387
    "B21" = 31.08109861050079; "C21" = 24.35828965023987; "D21" = 6.569594054913763; "E21" = 20.68537422931629; "G21" = 3.639805598281951; "N21" = 21.82391767718045
    "B22" = 31.89347353598117; "C22" = 25.11669099198368; "D22" = 6.675792623451391; "E22" = 21.4522705459396; "G22" = 3.633960465056922; "N22" = 21.95022141477427
    "B23" = 31.46146777671606; "C23" = 24.7140103772524; "D23" = 6.619112134825506; "E23" = 21.04595345041047; "G23" = 3.637064272775887; "N23" = 21.88276516421909
    "B24" = 29.79199376227679; "C24" = 23.14371967021903; "D24" = 6.404637674769112; "E24" = 19.44111215321573; "G24" = 3.649170313222658; "N24" = 21.62833583821499
    "B25" = 27.92896241018412; "C25" = 21.35908753257609; "D25" = 6.174938969342484; "E25" = 17.6115994931377; "G25" = 3.662994843394838; "N25" = 21.35644493247034
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}